$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.136.50"
Set-TextValue $ws.Range("E2") "  +0.06%  "
Set-TextValue $ws.Range("D3") "2.590.22"
Set-TextValue $ws.Range("E3") "  -1.16%  "
Set-TextValue $ws.Range("E4") "  +0.06%  "
Set-TextValue $ws.Range("D5") "519.52"
Set-TextValue $ws.Range("E5") "  +0.24%  "
Set-TextValue $ws.Range("D6") "144.01"
Set-TextValue $ws.Range("E6") "  +1.39%  "
Set-TextValue $ws.Range("E8") "  +0.08%  "
Set-TextValue $ws.Range("D9") "2.606.82"
Set-TextValue $ws.Range("E9") "  -0.66%  "
Set-TextValue $ws.Range("D10") "6.67"
Set-TextValue $ws.Range("E10") "  +0.29%  "
Set-TextValue $ws.Range("E11") "  -1.19%  "
Set-TextValue $ws.Range("E12") "  -2.94%  "
Set-TextValue $ws.Range("E13") "  -0.86%  "
Set-TextValue $ws.Range("D14") "3.047.17"
Set-TextValue $ws.Range("E14") "  -1.17%  "
Set-TextValue $ws.Range("D15") "58.078.84"
Set-TextValue $ws.Range("E15") "  -0.01%  "
Set-TextValue $ws.Range("D16") "20.44"
Set-TextValue $ws.Range("E16") "  -1.15%  "
Set-TextValue $ws.Range("B17") "WrappedEther"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "2.627.44"
Set-TextValue $ws.Range("E17") "  +0.18%  "
Set-TextValue $ws.Range("B18") "ShibaInu"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.0000134"
Set-TextValue $ws.Range("E18") "  -0.86%  "
Set-TextValue $ws.Range("D19") "340.71"
Set-TextValue $ws.Range("E19") "  +1.70%  "
Set-TextValue $ws.Range("D20") "4.33"
Set-TextValue $ws.Range("E20") "  -1.41%  "
Set-TextValue $ws.Range("D21") "10.30"
Set-TextValue $ws.Range("E21") "  -0.78%  "
Set-TextValue $ws.Range("D22") "6.37"
Set-TextValue $ws.Range("E22") "  +1.50%  "
Set-TextValue $ws.Range("D23") "0.998"
Set-TextValue $ws.Range("E23") "  -0.09%  "
Set-TextValue $ws.Range("D24") "66.22"
Set-TextValue $ws.Range("E24") "  +3.24%  "
Set-TextValue $ws.Range("E25") "  -0.69%  "
Set-TextValue $ws.Range("E26") "  -5.23%  "
Set-TextValue $ws.Range("D27") "0.998"
Set-TextValue $ws.Range("E27") "  -0.30%  "
Set-TextValue $ws.Range("D28") "2.706.36"
Set-TextValue $ws.Range("E28") "  -1.57%  "
Set-TextValue $ws.Range("D29") "7.02"
Set-TextValue $ws.Range("E29") "  -0.57%  "
Set-TextValue $ws.Range("D30") "0.0₃0751"
Set-TextValue $ws.Range("E30") "  -4.33%  "
Set-TextValue $ws.Range("E31") "  -0.02%  "
Set-TextValue $ws.Range("D32") "6.26"
Set-TextValue $ws.Range("E32") "  -5.10%  "
Set-TextValue $ws.Range("E33") "  +0.00%  "
Set-TextValue $ws.Range("D34") "18.77"
Set-TextValue $ws.Range("E34") "  +0.20%  "
Set-TextValue $ws.Range("D35") "149.83"
Set-TextValue $ws.Range("E35") "  -1.66%  "
Set-TextValue $ws.Range("E36") "  -1.47%  "
Set-TextValue $ws.Range("D37") "1.15"
Set-TextValue $ws.Range("E37") "  -1.96%  "
Set-TextValue $ws.Range("D38") "0.877"
Set-TextValue $ws.Range("E38") "  -2.89%  "
Set-TextValue $ws.Range("B39") "Stacks"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D39") "1.47"
Set-TextValue $ws.Range("E39") "  +2.12%  "
Set-TextValue $ws.Range("B40") "Fetch.AI"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D40") "0.840"
Set-TextValue $ws.Range("E40") "  -0.85%  "
Set-TextValue $ws.Range("D41") "36.09"
Set-TextValue $ws.Range("E41") "  -1.69%  "
Set-TextValue $ws.Range("E42") "  -1.30%  "
Set-TextValue $ws.Range("D43") "0.995"
Set-TextValue $ws.Range("E43") "  -0.43%  "
Set-TextValue $ws.Range("D44") "274.85"
Set-TextValue $ws.Range("E44") "  +2.27%  "
Set-TextValue $ws.Range("D45") "0.593"
Set-TextValue $ws.Range("E45") "  -0.95%  "
Set-TextValue $ws.Range("E46") "  +0.35%  "
Set-TextValue $ws.Range("D47") "0.0954"
Set-TextValue $ws.Range("E47") "  -1.37%  "
Set-TextValue $ws.Range("D48") "18.85"
Set-TextValue $ws.Range("E48") "  -2.13%  "
Set-TextValue $ws.Range("E49") "  -1.65%  "
Set-TextValue $ws.Range("D50") "4.70"
Set-TextValue $ws.Range("E50") "  +1.12%  "
Set-TextValue $ws.Range("D51") "1.980.75"
Set-TextValue $ws.Range("E51") "  -2.47%  "
